$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("company_list")

# Row 2
$ws.Range("D2").Value = 374568
$ws.Range("E2").Value = 5158
$ws.Range("F2").Value = 5158
$ws.Range("G2").Value = 372
$ws.Range("H2").Value = -1630
$ws.Range("I2").Value = -3645
$ws.Range("J2").Value = 2014
$ws.Range("K2").Value = 1236843
$ws.Range("L2").Value = 1115690
$ws.Range("M2").Value = 121152
$ws.Range("N2").Value = 43583
$ws.Range("O2").Value = 77569
$ws.Range("P2").Value = 3772
$ws.Range("Q2").Value = 45790
$ws.Range("R2").Value = -40350
$ws.Range("S2").Value = -6563
$ws.Range("T2").Value = 8214
$ws.Range("U2").Value = 37576
$ws.Range("V2").Value = 100622
$ws.Range("W2").Value = 1.38
$ws.Range("X2").Value = -0.43
$ws.Range("Y2").Value = -8.42
$ws.Range("Z2").Value = -0.14
$ws.Range("AA2").Value = 920.9
$ws.Range("AB2").Value = 888.09
$ws.Range("AC2").Value = -4831
$ws.Range("AD2").Value = -6.29
$ws.Range("AE2").Value = 62658
$ws.Range("AF2").Value = 0.49
$ws.Range("AG2").Value = 500
$ws.Range("AH2").Value = 1.64
$ws.Range("AI2").Value = -9.550000000000001
$ws.Range("AJ2").Value = 74958735

# Row 3
$ws.Range("D3").Value = 413763
$ws.Range("E3").Value = 7584
$ws.Range("F3").Value = 7584
$ws.Range("G3").Value = 1754
$ws.Range("H3").Value = 1205
$ws.Range("I3").Value = -2850
$ws.Range("J3").Value = 4055
$ws.Range("K3").Value = 1456217
$ws.Range("L3").Value = 1319760
$ws.Range("M3").Value = 136457
$ws.Range("N3").Value = 41545
$ws.Range("O3").Value = 94912
$ws.Range("P3").Value = 3772
$ws.Range("Q3").Value = 96207
$ws.Range("R3").Value = -89582
$ws.Range("S3").Value = 99
$ws.Range("T3").Value = 14621
$ws.Range("U3").Value = 81586
$ws.Range("V3").Value = 116822
$ws.Range("W3").Value = 1.83
$ws.Range("X3").Value = 0.29
$ws.Range("Y3").Value = -6.7
$ws.Range("Z3").Value = 0.09
$ws.Range("AA3").Value = 967.16
$ws.Range("AB3").Value = 824.9400000000001
$ws.Range("AC3").Value = -3778
$ws.Range("AD3").Value = -10.14
$ws.Range("AE3").Value = 59727
$ws.Range("AF3").Value = 0.64
$ws.Range("AG3").Value = 500
$ws.Range("AH3").Value = 1.3
$ws.Range("AI3").Value = -12.21
$ws.Range("AJ3").Value = 74958735

# Row 4
$ws.Range("D4").Value = 471202
$ws.Range("E4").Value = 16859
$ws.Range("F4").Value = 16859
$ws.Range("G4").Value = 16719
$ws.Range("H4").Value = 12887
$ws.Range("I4").Value = 4917
$ws.Range("J4").Value = 7970
$ws.Range("K4").Value = 1548710
$ws.Range("L4").Value = 1406431
$ws.Range("M4").Value = 142279
$ws.Range("N4").Value = 44081
$ws.Range("O4").Value = 98198
$ws.Range("P4").Value = 4896
$ws.Range("Q4").Value = 60973
$ws.Range("R4").Value = -58573
$ws.Range("S4").Value = 12947
$ws.Range("T4").Value = 11425
$ws.Range("U4").Value = 49548
$ws.Range("V4").Value = 118575
$ws.Range("W4").Value = 3.58
$ws.Range("X4").Value = 2.73
$ws.Range("Y4").Value = 11.48
$ws.Range("Z4").Value = 0.86
$ws.Range("AA4").Value = 988.5
$ws.Range("AB4").Value = 748.89
$ws.Range("AC4").Value = 6082
$ws.Range("AD4").Value = 5.76
$ws.Range("AE4").Value = 47899
$ws.Range("AF4").Value = 0.73
$ws.Range("AG4").Value = 600
$ws.Range("AH4").Value = 1.71
$ws.Range("AI4").Value = 11.6
$ws.Range("AJ4").Value = 74958735

# Row 5
$ws.Range("D5").Value = 504044
$ws.Range("E5").Value = 21589
$ws.Range("F5").Value = 21589
$ws.Range("G5").Value = 19951
$ws.Range("H5").Value = 13109
$ws.Range("I5").Value = 4054
$ws.Range("J5").Value = 9056
$ws.Range("K5").Value = 1601950
$ws.Range("L5").Value = 1442302
$ws.Range("M5").Value = 159647
$ws.Range("N5").Value = 45413
$ws.Range("O5").Value = 114235
$ws.Range("P5").Value = 4896
$ws.Range("Q5").Value = 51385
$ws.Range("R5").Value = -56658
$ws.Range("S5").Value = -3049
$ws.Range("T5").Value = 9678
$ws.Range("U5").Value = 41707
$ws.Range("V5").Value = 109442
$ws.Range("W5").Value = 4.28
$ws.Range("X5").Value = 2.6
$ws.Range("Y5").Value = 9.06
$ws.Range("Z5").Value = 0.83
$ws.Range("AA5").Value = 903.4299999999999
$ws.Range("AB5").Value = 799.15
$ws.Range("AC5").Value = 4140
$ws.Range("AD5").Value = 10.02
$ws.Range("AE5").Value = 49346
$ws.Range("AF5").Value = 0.84
$ws.Range("AG5").Value = 600
$ws.Range("AH5").Value = 1.45
$ws.Range("AI5").Value = 13.88
$ws.Range("AJ5").Value = 74958735

# Row 6
$ws.Range("D6").Value = 487402
$ws.Range("E6").Value = 18061
$ws.Range("F6").Value = 18061
$ws.Range("G6").Value = 12998
$ws.Range("H6").Value = 7993
$ws.Range("I6").Value = 4684
$ws.Range("K6").Value = 1695486
$ws.Range("L6").Value = 1525955
$ws.Range("M6").Value = 169530
$ws.Range("N6").Value = 41854
$ws.Range("P6").Value = 4896
$ws.Range("Q6").Value = 27492
$ws.Range("R6").Value = -26322
$ws.Range("S6").Value = 19913
$ws.Range("T6").Value = 15305
$ws.Range("U6").Value = 12187
$ws.Range("V6").Value = 128624
$ws.Range("W6").Value = 3.71
$ws.Range("X6").Value = 1.64
$ws.Range("Y6").Value = 10.73
$ws.Range("Z6").Value = 0.49
$ws.Range("AA6").Value = 900.11
$ws.Range("AB6").Value = 849.0700000000001
$ws.Range("AC6").Value = 4784
$ws.Range("AD6").Value = 6.55
$ws.Range("AE6").Value = 45478
$ws.Range("AF6").Value = 0.6899999999999999
$ws.Range("AG6").Value = 700
$ws.Range("AH6").Value = 2.23
$ws.Range("AI6").Value = 13.25
$ws.Range("AJ6").Value = 74958735

# Row 7
$ws.Range("D7").Value = 514983
$ws.Range("E7").Value = 13304
$ws.Range("G7").Value = 9500
$ws.Range("H7").Value = 7162
$ws.Range("I7").Value = 2998
$ws.Range("K7").Value = 1767535
$ws.Range("L7").Value = 1594280
$ws.Range("M7").Value = 181358
$ws.Range("N7").Value = 44908
$ws.Range("P7").Value = 4899
$ws.Range("Q7").Value = 27171
$ws.Range("R7").Value = -32070
$ws.Range("S7").Value = -1004
$ws.Range("T7").Value = 14119
$ws.Range("U7").Value = -16105
$ws.Range("W7").Value = 2.58
$ws.Range("X7").Value = 1.39
$ws.Range("Y7").Value = 6.91
$ws.Range("Z7").Value = 0.41
$ws.Range("AA7").Value = 879.08
$ws.Range("AC7").Value = 3063
$ws.Range("AD7").Value = 7.1
$ws.Range("AE7").Value = 48798
$ws.Range("AF7").Value = 0.45
$ws.Range("AG7").Value = 692
$ws.Range("AH7").Value = 3.18
$ws.Range("AI7").Value = 17.31

# Row 8
$ws.Range("D8").Value = 534036
$ws.Range("E8").Value = 16904
$ws.Range("G8").Value = 13115
$ws.Range("H8").Value = 9567
$ws.Range("I8").Value = 4252
$ws.Range("K8").Value = 1811355
$ws.Range("L8").Value = 1628526
$ws.Range("M8").Value = 191080
$ws.Range("N8").Value = 48245
$ws.Range("P8").Value = 4899
$ws.Range("Q8").Value = 33527
$ws.Range("R8").Value = -29342
$ws.Range("S8").Value = -2157
$ws.Range("T8").Value = 11273
$ws.Range("U8").Value = 224
$ws.Range("W8").Value = 3.16
$ws.Range("X8").Value = 1.79
$ws.Range("Y8").Value = 9.130000000000001
$ws.Range("Z8").Value = 0.53
$ws.Range("AA8").Value = 852.28
$ws.Range("AC8").Value = 4342
$ws.Range("AD8").Value = 5.01
$ws.Range("AE8").Value = 52423
$ws.Range("AF8").Value = 0.41
$ws.Range("AG8").Value = 701
$ws.Range("AH8").Value = 3.22
$ws.Range("AI8").Value = 12.36

# Row 9
$ws.Range("D9").Value = 553960
$ws.Range("E9").Value = 19193
$ws.Range("G9").Value = 15403
$ws.Range("H9").Value = 11218
$ws.Range("I9").Value = 5336
$ws.Range("K9").Value = 1842057
$ws.Range("L9").Value = 1655090
$ws.Range("M9").Value = 203723
$ws.Range("N9").Value = 52197
$ws.Range("P9").Value = 4899
$ws.Range("Q9").Value = 38869
$ws.Range("R9").Value = -24156
$ws.Range("S9").Value = -3376
$ws.Range("T9").Value = 11372
$ws.Range("U9").Value = 2181
$ws.Range("W9").Value = 3.46
$ws.Range("X9").Value = 2.02
$ws.Range("Y9").Value = 10.63
$ws.Range("Z9").Value = 0.61
$ws.Range("AA9").Value = 812.42
$ws.Range("AC9").Value = 5450
$ws.Range("AD9").Value = 3.99
$ws.Range("AE9").Value = 56718
$ws.Range("AF9").Value = 0.38
$ws.Range("AG9").Value = 723
$ws.Range("AH9").Value = 3.32
$ws.Range("AI9").Value = 10.15
